# merge database Quan Ly Nhan Su into ISO
# Adds two new rows (DAO layer entries) to Sheet1's function list, and
# refreshes the workbook's default font/formatting to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Default font: Arial -> Calibri (affects every cell via the Normal
#    style, which is what every cellXf in this workbook points to).
# ---------------------------------------------------------------------
$wb.Styles.Item(1).Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 2. Column widths shifted slightly (consequence of the font metrics
#    change recalculating the stored character widths).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 4.711495535714286
$ws.Columns.Item(3).ColumnWidth = 14.570870535714286
$ws.Columns.Item(4).ColumnWidth = 27.711495535714285
$ws.Columns.Item(5).ColumnWidth = 31.570870535714285
$ws.Columns.Item(6).ColumnWidth = 61.141183035714285
$ws.Columns.Item(7).ColumnWidth = 11.711495535714286

# ---------------------------------------------------------------------
# 3. A handful of existing rows grew slightly taller.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 135
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 120
$ws.Rows.Item(8).RowHeight = 135

# ---------------------------------------------------------------------
# 4. New rows 10 & 11: two more DAO entries merged in from the
#    "Quan Ly Nhan Su" workbook.
#
#    Cells are written in a specific order so that the new shared
#    strings land at the same indices the source workbook used:
#      22 getKhoaByMaBoPhan()   23 DAO   24 ThanhVienDAO
#      25 UpdateVaiTroTV()      26 KhoaDAO
# ---------------------------------------------------------------------
$ws.Range("F10").Value = "getKhoaByMaBoPhan()"
$ws.Range("D10").Value = "DAO"
$ws.Range("E11").Value = "ThanhVienDAO"
$ws.Range("F11").Value = "UpdateVaiTroTV()"
$ws.Range("E10").Value = "KhoaDAO"

$ws.Range("C10").Value = 8
$ws.Range("G10").Value = "Thêm"
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = "DAO"
$ws.Range("G11").Value = "Thêm"

# Match the formatting pattern already used by the rows above
# (vertical-top on C/D, wrap-text on E, plain on F, vertical-top on G).
$ws.Range("C10:D10").VerticalAlignment = -4160
$ws.Range("E10").WrapText = $true
$ws.Range("C11:D11").VerticalAlignment = -4160
$ws.Range("E11").WrapText = $true
